$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = "Variação 2022/2013"
$ws.Range("C10").Value = 7.849957117115469
